# Correction of English text on the "Fuel" column (column E) of Sheet1:
#   "Super"    -> "Premium"
#   "Ordinary" -> "Regular"
#   "Diesel" stays "Diesel" (unchanged)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("Super", "Premium")
$ws.Cells.Replace("Ordinary", "Regular")

# Reflect the author's final selection (column E was selected after the edit).
$ws.Columns("E").Select()
